$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.630832113316842
$ws.Cells.Item(2, 4).Value = 0.02880582012058142
$ws.Cells.Item(2, 5).Value = 0.1655006580538818
$ws.Cells.Item(2, 6).Value = 0.9216479529410009
$ws.Cells.Item(2, 7).Value = 0.7730921049688391
$ws.Cells.Item(2, 8).Value = 0.8311006263541998
$ws.Cells.Item(2, 9).Value = 0.9894943373427445
$ws.Cells.Item(2, 11).Value = 0.5375561185338711
$ws.Cells.Item(2, 12).Value = 0.203407482233942
$ws.Cells.Item(2, 14).Value = 1.755416247555921
$ws.Cells.Item(3, 2).Value = 0.6065567131176977
$ws.Cells.Item(3, 4).Value = 0.02807919792270042
$ws.Cells.Item(3, 5).Value = 0.1661664880952003
$ws.Cells.Item(3, 6).Value = 0.9119622753842123
$ws.Cells.Item(3, 7).Value = 0.7644004968154121
$ws.Cells.Item(3, 8).Value = 0.8317218170034124
$ws.Cells.Item(3, 9).Value = 0.9966014868725495
$ws.Cells.Item(3, 11).Value = 0.4701167810797813
$ws.Cells.Item(3, 12).Value = 0.1909835505091166
$ws.Cells.Item(3, 14).Value = 1.774230150580513
$ws.Cells.Item(4, 2).Value = 0.5919269093389232
$ws.Cells.Item(4, 4).Value = 0.02762798263261956
$ws.Cells.Item(4, 5).Value = 0.1666046963100674
$ws.Cells.Item(4, 6).Value = 0.9066079626795229
$ws.Cells.Item(4, 7).Value = 0.759605367938164
$ws.Cells.Item(4, 8).Value = 0.8325224494669214
$ws.Cells.Item(4, 9).Value = 1.001402809579425
$ws.Cells.Item(4, 11).Value = 0.4286897552950393
$ws.Cells.Item(4, 12).Value = 0.1834602331612274
$ws.Cells.Item(4, 14).Value = 1.786371471599869
$ws.Cells.Item(5, 2).Value = 0.5860348481664062
$ws.Cells.Item(5, 4).Value = 0.02744284522351847
$ws.Cells.Item(5, 5).Value = 0.1667906762306595
$ws.Cells.Item(5, 6).Value = 0.9045749419929479
$ws.Cells.Item(5, 7).Value = 0.7577872463727573
$ws.Cells.Item(5, 8).Value = 0.8329541244520868
$ws.Cells.Item(5, 9).Value = 1.003469415697367
$ws.Cells.Item(5, 11).Value = 0.4118034328219835
$ws.Cells.Item(5, 12).Value = 0.1804208777293184
$ws.Cells.Item(5, 14).Value = 1.791467344747128
$ws.Cells.Item(6, 2).Value = 0.5850606991334928
$ws.Cells.Item(6, 4).Value = 0.02741202729313486
$ws.Cells.Item(6, 5).Value = 0.1668220058979855
$ws.Cells.Item(6, 6).Value = 0.9042463498408182
$ws.Cells.Item(6, 7).Value = 0.7574935514909384
$ws.Cells.Item(6, 8).Value = 0.8330321701810135
$ws.Cells.Item(6, 9).Value = 1.003819219926868
$ws.Cells.Item(6, 11).Value = 0.408999210731821
$ws.Cells.Item(6, 12).Value = 0.179917794459385
$ws.Cells.Item(6, 14).Value = 1.792322459253066
$ws.Cells.Item(7, 2).Value = 0.5918471641157055
$ws.Cells.Item(7, 4).Value = 0.02762549090699906
$ws.Cells.Item(7, 5).Value = 0.1666071744878943
$ws.Cells.Item(7, 6).Value = 0.9065799419303886
$ws.Cells.Item(7, 7).Value = 0.7595802980571591
$ws.Cells.Item(7, 8).Value = 0.8325278443919046
$ws.Cells.Item(7, 9).Value = 1.001430235034995
$ws.Cells.Item(7, 11).Value = 0.4284620381503998
$ws.Cells.Item(7, 12).Value = 0.1834191361617599
$ws.Cells.Item(7, 14).Value = 1.786439596396033
$ws.Cells.Item(8, 2).Value = 0.6224050611502037
$ws.Cells.Item(8, 4).Value = 0.02855633834754912
$ws.Cells.Item(8, 5).Value = 0.1657241475261813
$ws.Cells.Item(8, 6).Value = 0.91818521828813
$ws.Cells.Item(8, 7).Value = 0.7699826860930159
$ws.Cells.Item(8, 8).Value = 0.8312277985958048
$ws.Cells.Item(8, 9).Value = 0.9918540866977672
$ws.Cells.Item(8, 11).Value = 0.5143071406344006
$ws.Cells.Item(8, 12).Value = 0.1991019365206483
$ws.Cells.Item(8, 14).Value = 1.76178090023776
$ws.Cells.Item(9, 2).Value = 0.6844977810303021
$ws.Cells.Item(9, 4).Value = 0.03034114917130282
$ws.Cells.Item(9, 5).Value = 0.1642249348539486
$ws.Cells.Item(9, 6).Value = 0.9456551895619754
$ws.Cells.Item(9, 7).Value = 0.7946920479179624
$ws.Cells.Item(9, 8).Value = 0.8320059657979755
$ws.Cells.Item(9, 9).Value = 0.9765462843245523
$ws.Cells.Item(9, 11).Value = 0.6824933106811386
$ws.Cells.Item(9, 12).Value = 0.2306889013077011
$ws.Cells.Item(9, 14).Value = 1.71810425943672
$ws.Cells.Item(10, 2).Value = 0.7314228318133473
$ws.Cells.Item(10, 4).Value = 0.03162732032737381
$ws.Cells.Item(10, 5).Value = 0.1632641033629265
$ws.Cells.Item(10, 6).Value = 0.9687261447779605
$ws.Cells.Item(10, 7).Value = 0.815495373888794
$ws.Cells.Item(10, 8).Value = 0.8346094163364484
$ws.Cells.Item(10, 9).Value = 0.9674154380749727
$ws.Cells.Item(10, 11).Value = 0.8059702099044728
$ws.Cells.Item(10, 12).Value = 0.2544059331004433
$ws.Cells.Item(10, 14).Value = 1.688868435910308
$ws.Cells.Item(11, 2).Value = 0.7530505451268255
$ws.Cells.Item(11, 4).Value = 0.03220690152753747
$ws.Cells.Item(11, 5).Value = 0.1628573174807286
$ws.Cells.Item(11, 6).Value = 0.9798527036754763
$ws.Cells.Item(11, 7).Value = 0.8255394951164305
$ws.Cells.Item(11, 8).Value = 0.8362357601100712
$ws.Cells.Item(11, 9).Value = 0.9637210042620339
$ws.Cells.Item(11, 11).Value = 0.862125545622888
$ws.Cells.Item(11, 12).Value = 0.265306824510688
$ws.Cells.Item(11, 14).Value = 1.676188109702954
$ws.Cells.Item(12, 2).Value = 0.7612804592449436
$ws.Cells.Item(12, 4).Value = 0.03242557349994257
$ws.Cells.Item(12, 5).Value = 0.1627076187418033
$ws.Cells.Item(12, 6).Value = 0.984157071100384
$ws.Cells.Item(12, 7).Value = 0.8294267574007108
$ws.Cells.Item(12, 8).Value = 0.8369152175614829
$ws.Cells.Item(12, 9).Value = 0.9623880573909531
$ws.Cells.Item(12, 11).Value = 0.8833878349480813
$ws.Cells.Item(12, 12).Value = 0.269450790639695
$ws.Cells.Item(12, 14).Value = 1.671475503208757
$ws.Cells.Item(13, 2).Value = 0.7595062299875508
$ws.Cells.Item(13, 4).Value = 0.03237851446415618
$ws.Cells.Item(13, 5).Value = 0.1627396661586151
$ws.Cells.Item(13, 6).Value = 0.9832259993342944
$ws.Cells.Item(13, 7).Value = 0.8285858356259297
$ws.Cells.Item(13, 8).Value = 0.8367660552221565
$ws.Cells.Item(13, 9).Value = 0.9626721932458082
$ws.Cells.Item(13, 11).Value = 0.8788087383940422
$ws.Cells.Item(13, 12).Value = 0.2685576011047459
$ws.Cells.Item(13, 14).Value = 1.672486480102876
$ws.Cells.Item(14, 2).Value = 0.7537268268240496
$ws.Cells.Item(14, 4).Value = 0.03222490795632638
$ws.Cells.Item(14, 5).Value = 0.1628449147356803
$ws.Cells.Item(14, 6).Value = 0.9802050021644249
$ws.Cells.Item(14, 7).Value = 0.8258576219890728
$ws.Cells.Item(14, 8).Value = 0.8362903847183958
$ws.Cells.Item(14, 9).Value = 0.9636100179761797
$ws.Cells.Item(14, 11).Value = 0.8638748595455752
$ws.Cells.Item(14, 12).Value = 0.265647430047693
$ws.Cells.Item(14, 14).Value = 1.675798613187064
$ws.Cells.Item(15, 2).Value = 0.7501919664011609
$ws.Cells.Item(15, 4).Value = 0.03213071464705308
$ws.Cells.Item(15, 5).Value = 0.1629099475741476
$ws.Cells.Item(15, 6).Value = 0.9783664078402836
$ws.Cells.Item(15, 7).Value = 0.8241974303278283
$ws.Cells.Item(15, 8).Value = 0.8360073059542117
$ws.Cells.Item(15, 9).Value = 0.9641930655526281
$ws.Cells.Item(15, 11).Value = 0.854727088657711
$ws.Cells.Item(15, 12).Value = 0.2638669528493267
$ws.Cells.Item(15, 14).Value = 1.677839004878248
$ws.Cells.Item(16, 2).Value = 0.7300150289187854
$ws.Cells.Item(16, 4).Value = 0.03158933159957655
$ws.Cells.Item(16, 5).Value = 0.1632912962254403
$ws.Cells.Item(16, 6).Value = 0.9680117198264639
$ws.Cells.Item(16, 7).Value = 0.8148506740816117
$ws.Cells.Item(16, 8).Value = 0.8345120284513428
$ws.Cells.Item(16, 9).Value = 0.9676661204808283
$ws.Cells.Item(16, 11).Value = 0.802299994575236
$ws.Cells.Item(16, 12).Value = 0.2536957811087461
$ws.Cells.Item(16, 14).Value = 1.689709594893786
$ws.Cells.Item(17, 2).Value = 0.7177088315311835
$ws.Cells.Item(17, 4).Value = 0.03125579272265355
$ws.Cells.Item(17, 5).Value = 0.1635329913781212
$ws.Cells.Item(17, 6).Value = 0.9618213136222238
$ws.Cells.Item(17, 7).Value = 0.8092656456409912
$ws.Cells.Item(17, 8).Value = 0.8337079604722533
$ws.Cells.Item(17, 9).Value = 0.9699143605012637
$ws.Cells.Item(17, 11).Value = 0.7701335835802467
$ws.Cells.Item(17, 12).Value = 0.2474847164564693
$ws.Cells.Item(17, 14).Value = 1.697150491542398
$ws.Cells.Item(18, 2).Value = 0.7106571355183462
$ws.Cells.Item(18, 4).Value = 0.03106343254120247
$ws.Cells.Item(18, 5).Value = 0.1636748610170324
$ws.Cells.Item(18, 6).Value = 0.9583201816510467
$ws.Cells.Item(18, 7).Value = 0.8061079153073081
$ws.Cells.Item(18, 8).Value = 0.8332870848432066
$ws.Cells.Item(18, 9).Value = 0.9712507128724539
$ws.Cells.Item(18, 11).Value = 0.7516309321389656
$ws.Cells.Item(18, 12).Value = 0.2439228117282539
$ws.Cells.Item(18, 14).Value = 1.701488572925239
$ws.Cells.Item(19, 2).Value = 0.7082741189466617
$ws.Cells.Item(19, 4).Value = 0.03099821416405035
$ws.Cells.Item(19, 5).Value = 0.163723386101345
$ws.Cells.Item(19, 6).Value = 0.9571449595948138
$ws.Cells.Item(19, 7).Value = 0.8050481362204209
$ws.Cells.Item(19, 8).Value = 0.8331517279158476
$ws.Cells.Item(19, 9).Value = 0.9717106020868016
$ws.Cells.Item(19, 11).Value = 0.7453660261288064
$ws.Cells.Item(19, 12).Value = 0.2427186231262795
$ws.Cells.Item(19, 14).Value = 1.702967377890776
$ws.Cells.Item(20, 2).Value = 0.719016107731818
$ws.Cells.Item(20, 4).Value = 0.03129135214102519
$ws.Cells.Item(20, 5).Value = 0.1635069673539427
$ws.Cells.Item(20, 6).Value = 0.9624741412830673
$ws.Cells.Item(20, 7).Value = 0.809854525360123
$ws.Cells.Item(20, 8).Value = 0.8337892488777783
$ws.Cells.Item(20, 9).Value = 0.9696705578441396
$ws.Cells.Item(20, 11).Value = 0.7735578987303882
$ws.Cells.Item(20, 12).Value = 0.2481448049185246
$ws.Cells.Item(20, 14).Value = 1.69635236335591
$ws.Cells.Item(21, 2).Value = 0.7554232967033272
$ws.Cells.Item(21, 4).Value = 0.03227004778968023
$ws.Cells.Item(21, 5).Value = 0.1628138829489392
$ws.Cells.Item(21, 6).Value = 0.9810898716528982
$ws.Cells.Item(21, 7).Value = 0.8266566886652242
$ws.Cells.Item(21, 8).Value = 0.836428374411966
$ws.Cells.Item(21, 9).Value = 0.9633327632611355
$ws.Cells.Item(21, 11).Value = 0.868261373052178
$ws.Cells.Item(21, 12).Value = 0.2665017830837257
$ws.Cells.Item(21, 14).Value = 1.674823338737459
$ws.Cells.Item(22, 2).Value = 0.779450244596859
$ws.Cells.Item(22, 4).Value = 0.03290499717726902
$ws.Cells.Item(22, 5).Value = 0.1623862153184956
$ws.Cells.Item(22, 6).Value = 0.9937867169059444
$ws.Cells.Item(22, 7).Value = 0.8381262903645705
$ws.Cells.Item(22, 8).Value = 0.8385239002459457
$ws.Cells.Item(22, 9).Value = 0.9595756759890079
$ws.Cells.Item(22, 11).Value = 0.9301407213837081
$ws.Cells.Item(22, 12).Value = 0.2785925686382029
$ws.Cells.Item(22, 14).Value = 1.6612726137494
$ws.Cells.Item(23, 2).Value = 0.7666054719719568
$ws.Cells.Item(23, 4).Value = 0.03256654510997237
$ws.Cells.Item(23, 5).Value = 0.1626121592615921
$ws.Cells.Item(23, 6).Value = 0.9869615824943452
$ws.Cells.Item(23, 7).Value = 0.8319599639976758
$ws.Cells.Item(23, 8).Value = 0.8373715476368915
$ws.Cells.Item(23, 9).Value = 0.9615456681142263
$ws.Cells.Item(23, 11).Value = 0.897116029803442
$ws.Cells.Item(23, 12).Value = 0.272130957481977
$ws.Cells.Item(23, 14).Value = 1.668457294887219
$ws.Cells.Item(24, 2).Value = 0.7184250155467566
$ws.Cells.Item(24, 4).Value = 0.0312752776097156
$ws.Cells.Item(24, 5).Value = 0.1635187237327964
$ws.Cells.Item(24, 6).Value = 0.9621788178183976
$ws.Cells.Item(24, 7).Value = 0.809588127215946
$ws.Cells.Item(24, 8).Value = 0.8337523694485043
$ws.Cells.Item(24, 9).Value = 0.9697806445613608
$ws.Cells.Item(24, 11).Value = 0.772009796209062
$ws.Cells.Item(24, 12).Value = 0.247846351140268
$ws.Cells.Item(24, 14).Value = 1.69671300955754
$ws.Cells.Item(25, 2).Value = 0.667469592660126
$ws.Cells.Item(25, 4).Value = 0.02986269560136634
$ws.Cells.Item(25, 5).Value = 0.1646057394712992
$ws.Cells.Item(25, 6).Value = 0.9377177053686125
$ws.Cells.Item(25, 7).Value = 0.7875438223978364
$ws.Cells.Item(25, 8).Value = 0.8314388896689877
$ws.Cells.Item(25, 9).Value = 0.9803158810099823
$ws.Cells.Item(25, 11).Value = 0.637010607746987
$ws.Cells.Item(25, 12).Value = 0.2220543745433474
$ws.Cells.Item(25, 14).Value = 1.72941894248982
